$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table (rows 2-51) with the
# latest scraped values. Several "price" cells contain values that look
# numeric (e.g. "1.000", "307.45") but must remain plain text, matching
# the original inlineStr cell contents. We force text by prefixing the
# value with a leading apostrophe (Excel's "treat as text" marker) and
# then reset the cell style back to Normal so no stray formatting/
# quote-prefix styling is left behind.

$ws.Range("D2").Value = "27.094.20"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "1.887.93"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'307.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5151"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("D8").Value = "'0.3718"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "'0.07212"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").Value = "'0.9035"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "'21.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").Value = "'0.07618"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "1.883.90"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "'94.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").Value = "'5.273"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'0.000008501"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").Value = "'0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "27.136.01"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "'5.053"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").Value = "2.139.50"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'145.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "'1.792"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'18.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").Value = "'2.173"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.82%  "
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").Value = "'4.978"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.07%  "
$ws.Range("D31").Value = "'4.823"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'0.09208"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'0.05067"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'1.198"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.86%  "
$ws.Range("D35").Value = "'0.7588"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").Value = "'3.003"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").Value = "'3.272"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "'2.564"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("D39").Value = "'0.5624"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.79%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").Value = "'1.076"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "'9.029"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.15%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'118.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.574"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("E45").Value = "  +3.71%  "
$ws.Range("D46").Value = "'0.4800"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'0.9994"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'1.576"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'63.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.26%  "
